# Apply RQ5 variable-definitions edit:
# Rename "Delegator" terminology to "Participant" in most places, and
# repoint the "Source" row away from the delegation survey toward the
# text-scenario source (per commit message: "updated RQ5 and variables
# to base off text instead of delegation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (column headers) ---
$ws.Range("B1").Value = "Participant_ID"
$ws.Range("H1").Value = "Alignment score (Participant|Most aligned target)"
$ws.Range("I1").Value = "Alignment score (Participant|Least aligned target)"
$ws.Range("K1").Value = "Alignment score (Participant|group target)"
$ws.Range("M1").Value = "Alignment score (Participant|ADM(most))"
$ws.Range("N1").Value = "Alignment score (Participant|ADM(least))"

# --- Row 2 (Source) ---
$ws.Range("B2").Value = "Text scenario"
$ws.Range("C2").Value = "Text scenario"
$ws.Range("D2").Value = "Text scenario"
$ws.Range("E2").Value = "Text scenario loading log"
# TA2_Name no longer has a Source entry
$ws.Range("L2").Clear()

# --- Row 3 (Definition) ---
$ws.Range("B3").Value = "Used to track and identify participants, also called delegator ID"
$ws.Range("E3").Value = "Scenario presented to the participant in the text scenarios"
$ws.Range("F3").Value = "Target with the highest alignment score for the participant on the text scenario"
$ws.Range("G3").Value = "Target with the lowest alignment score for the participant on the text scenario"
$ws.Range("H3").Value = "Calculated alignment score between the participant and a target"
$ws.Range("I3").Value = "Calculated alignment score between the participant and a target"
$ws.Range("K3").Value = "Calculated alignment score between the KDMA measurement of a participant and a group target"
# New definition for TA2_Name
$ws.Range("L3").Value = "Source of ADM responses"
$ws.Range("M3").Value = "Calculated alignment score between the participant and the aligned ADM run on the same scenario at the most aligned target"
$ws.Range("N3").Value = "Calculated alignment score between the participant and the aligned ADM run on the same scenario at the least aligned target"
$ws.Range("O3").Value = "% of exact matches on probe responses between participant and ADM run on same scenario at most aligned target"
$ws.Range("P3").Value = "% of exact matches on probe responses between participant and ADM run on same scenario at least aligned target"
$ws.Range("Q3").Value = "% of exact matches on probe responses among group members and ADM run on same scenario at group target"

# --- Sheet view: scrolled right a bit, selection moved to N2 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$ws.Range("N2").Select()
